$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Rondônia) - date + value update only
$ws.Range("C2").Value = "'01/07/2024"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = 97.8211009174312

# Row 3 - region swapped to Mato Grosso, date + value update
$ws.Range("A3").Value = "Mato Grosso"
$ws.Range("C3").Value = "'01/07/2024"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = 97.72382397572079

# Row 4 - region swapped to Santa Catarina, date + value update
$ws.Range("A4").Value = "Santa Catarina"
$ws.Range("C4").Value = "'01/07/2024"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = 97.23912026204961

# Row 5 (Mato Grosso do Sul) - date + value update only
$ws.Range("C5").Value = "'01/07/2024"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").Value = 96.53102068045364

# Row 6 - region swapped to Paraná, date + value update
$ws.Range("A6").Value = "Paraná"
$ws.Range("C6").Value = "'01/07/2024"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = 95.98285169895205

# Row 7 - region swapped to Espírito Santo, date + value update
$ws.Range("A7").Value = "Espírito Santo"
$ws.Range("C7").Value = "'01/07/2024"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").Value = 95.8891454965358

# Row 8 (Sergipe) - date + value + rank update
$ws.Range("C8").Value = "'01/07/2024"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = 91.61462979482606
$ws.Range("E8").Value = "22º"

# Row 9 (Nordeste) - date + value update only
$ws.Range("C9").Value = "'01/07/2024"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = 91.32792385113521

# Row 10 (Brasil) - date + value update only
$ws.Range("C10").Value = "'01/07/2024"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = 93.63718985731164
